# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Adds a new "Number of employees / Assets / Turnover" breakdown table
# (by enterprise size class: Micro / Small / Medium / Large) to the
# Iceland MSME summary sheet, inserted just above the existing
# "SME Performance Review EU" source block, which shifts down from
# rows 26:27 to rows 32:33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing "SME Performance Review EU" source block (rows 26:27)
# down to rows 32:33, opening up six blank rows (23:28) for the new table.
$ws.Range("23:28").Insert()

# --- Row 23: new table header (bold, matching the "title" look used by the
#             B11:D11 / B19:D19 header rows above) -------------------------
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Font.Bold = $true

# --- Row 24: Micro ----------------------------------------------------------
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

# --- Row 25: Small -----------------------------------------------------------
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

# --- Row 26: Medium -----------------------------------------------------------
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

# --- Row 27: Large ------------------------------------------------------------
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
